$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cells (B16:B22) hold numeric-looking values that are stored
# as TEXT (shared strings), e.g. "42381". A plain .Value assignment would
# have Excel auto-detect the new text as a number and convert the cell to
# a numeric cell, which would not match the desired edit (the cells must
# stay text). Formatting the range as Text ("@") first makes Excel keep
# replaced numeric-looking text as text, exactly like typing into a
# Text-formatted cell in the real application.
$ws.Range("B16:B22").NumberFormat = "@"

# Replace each old text value with its new value (old -> new), matching
# the shared-string content changes from the diff, row by row.
$ws.Cells.Replace("42381", "43170")
$ws.Cells.Replace("42991", "43617")
$ws.Cells.Replace("42424", "42150")
$ws.Cells.Replace("42768", "43119")
$ws.Cells.Replace("43124", "42976")
$ws.Cells.Replace("42567", "41750")
$ws.Cells.Replace("41820", "42512")

# Restore the cells' style to the workbook default so no visible
# formatting change is left behind (cells go back to being unstyled,
# same as before, while keeping their new text values).
$ws.Range("B16:B22").Style = "Normal"
